# JCF_TB_1794_Heterogene_Begriffe.xlsx
# Fix the typo in the report title cell: "Heterogene Begrife" -> "Heterogene Begriffe"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Heterogene Begriffe"
